$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first pass of the permutation chain ---
# Sai_Nome / Entra_Nome values shift so that Marcos enters where Joao used to be.
$ws.Range("B2").Value = "Marcos"
$ws.Range("C2").Value = "João"

# --- Row 3 (new): second pass of the chain, continuing the swap ---
# Use a text formula first so the date-looking string "2025-12-01" is not
# auto-converted into a date serial value, then flatten it to a plain value.
$ws.Range("A3").Formula = "=""2025-12-01"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("B3").Value = "João"
$ws.Range("C3").Value = "Felipe"

# Copy the row 2 formatting (style + row height) onto the new row 3.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

$excel.CutCopyMode = 0
